# Updates cryptos list prices / volume(1h) figures, and reorders a couple
# of coin rows, matching the latest scrape from coinranking.com.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.257.83"
$ws.Range("E2").Value = "  +1.74%  "
$ws.Range("D3").Value = "3.226.90"
$ws.Range("E3").Value = "  +1.45%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'606.09"
$ws.Range("E5").Value = "  +4.75%  "
$ws.Range("D6").Value = "'155.31"
$ws.Range("E6").Value = "  +2.85%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "3.225.22"
$ws.Range("E8").Value = "  +1.44%  "
$ws.Range("D9").Value = "'0.537"
$ws.Range("E9").Value = "  +0.89%  "
$ws.Range("E10").Value = "  -0.90%  "
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("D13").Value = "'0.0000273"
$ws.Range("E13").Value = "  -0.55%  "
$ws.Range("D14").Value = "'39.10"
$ws.Range("E14").Value = "  +2.52%  "
$ws.Range("D15").Value = "3.753.16"
$ws.Range("E15").Value = "  +1.36%  "
$ws.Range("D16").Value = "'7.52"
$ws.Range("E16").Value = "  +4.63%  "
$ws.Range("D17").Value = "66.325.92"
$ws.Range("E17").Value = "  +1.73%  "
$ws.Range("D18").Value = "3.226.07"
$ws.Range("E18").Value = "  +1.21%  "
$ws.Range("D20").Value = "'515.11"
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("E21").Value = "  +6.30%  "
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("D23").Value = "'15.30"
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("D24").Value = "'8.03"
$ws.Range("E24").Value = "  +2.41%  "
$ws.Range("D25").Value = "'85.77"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").Value = "  +3.94%  "
$ws.Range("D28").Value = "'9.30"
$ws.Range("E28").Value = "  +2.62%  "
$ws.Range("E29").Value = "  +2.78%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "'6.93"
$ws.Range("E30").Value = "  +10.68%  "
$ws.Range("B31").Value = "Stacks"
$ws.Range("C31").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D31").Value = "'2.89"
$ws.Range("E31").Value = "  +3.70%  "
$ws.Range("D32").Value = "'28.36"
$ws.Range("E32").Value = "  +0.94%  "
$ws.Range("D33").Value = "'1.23"
$ws.Range("E33").Value = "  +1.47%  "
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").Value = "'6.69"
$ws.Range("E35").Value = "  +0.75%  "
$ws.Range("D36").Value = "'55.68"
$ws.Range("D37").Value = "'0.0925"
$ws.Range("E37").Value = "  +2.14%  "
$ws.Range("D38").Value = "'491.51"
$ws.Range("E38").Value = "  +2.67%  "
$ws.Range("D39").Value = "'0.0425"
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("D40").Value = "'3.05"
$ws.Range("E40").Value = "  -2.66%  "
$ws.Range("D41").Value = "'8.90"
$ws.Range("E41").Value = "  +2.64%  "
$ws.Range("D42").Value = "'0.298"
$ws.Range("E42").Value = "  +3.14%  "
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").Value = "'2.52"
$ws.Range("E44").Value = "  +4.42%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.997.75"
$ws.Range("E45").Value = "  -2.44%  "
$ws.Range("D46").Value = "0.0₃0652"
$ws.Range("E46").Value = "  +7.66%  "
$ws.Range("D47").Value = "'29.31"
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("E50").Value = "  +3.19%  "
$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").Value = "'33.83"
$ws.Range("E51").Value = "  +4.72%  "
